$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: strip the "27 - " wilaya-code prefix from the wilaya name.
$ws.Range("C19").Value = "Mostaganem"

# Append four new data rows (20-23), extending the used range to A1:E23.
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 23
$ws.Range("C20").Value = "Constantine"
$ws.Range("D20").Value = "Constantine"
$ws.Range("E20").Value = "6663b990f18db5d19ddbc69b"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 24
$ws.Range("C21").Value = "Constantine"
$ws.Range("D21").Value = "Constantine"
$ws.Range("E21").Value = "6663ba09f18db5d19ddbc69f"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 24
$ws.Range("C22").Value = "Constantine"
$ws.Range("D22").Value = "Constantine"
$ws.Range("E22").Value = "6663c158114e26a841e7b707"

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 24
$ws.Range("C23").Value = "25 - Constantine"
$ws.Range("D23").Value = "Constantine"
$ws.Range("E23").Value = "6663c30a114e26a841e7b86f"
